$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 511.44446
$ws.Range("I98").Value = 515
$ws.Range("K98").Value = 515
$ws.Range("M98").Value = 983
$ws.Range("H100").Value = 1999.7142
$ws.Range("I100").Value = 2136.5789
$ws.Range("J100").Value = 699.5
$ws.Range("K100").Value = 2136.5789
$ws.Range("L100").Value = 699.5
$ws.Range("M100").Value = -1595.5789
$ws.Range("N100").Value = -1781.5
$ws.Range("H103").Value = 4026.0435
$ws.Range("J103").Value = 5999.8887
$ws.Range("L103").Value = 17999.6661
$ws.Range("N103").Value = -19171.6661
$ws.Range("H122").Value = 511.44446
$ws.Range("I122").Value = 515
$ws.Range("K122").Value = 1545
$ws.Range("M122").Value = 905
$ws.Range("H137").Value = 3347.1875
$ws.Range("I137").Value = 1648.5
$ws.Range("K137").Value = 4945.5
$ws.Range("M137").Value = -2395.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 203596.25
$ws.Range("I6").Value = 571703.5600000001
$ws.Range("J6").Value = 5384.615
$ws.Range("K6").Value = 571703.5600000001
$ws.Range("L6").Value = 5384.615
$ws.Range("M6").Value = -571530.5600000001
$ws.Range("N6").Value = -5730.615
$ws.Range("H8").Value = 4012491
$ws.Range("I8").Value = 6683335
$ws.Range("K8").Value = 6683335
$ws.Range("M8").Value = -6683191
$ws.Range("H74").Value = 1681.3529
$ws.Range("I74").Value = 1648.0769
$ws.Range("K74").Value = 1648.0769
$ws.Range("M74").Value = -774.0769
$ws.Range("H77").Value = 1681.3529
$ws.Range("I77").Value = 1648.0769
$ws.Range("K77").Value = 8240.3845
$ws.Range("M77").Value = -3872.3845
$ws.Range("H92").Value = 52000.25
$ws.Range("J92").Value = 49333.668
$ws.Range("L92").Value = 49333.668
$ws.Range("N92").Value = -54325.668
$ws.Range("H110").Value = 1511
$ws.Range("I110").Value = 1590.125
$ws.Range("K110").Value = 1590.125
$ws.Range("M110").Value = 454.875
$ws.Range("H132").Value = 3853.2727
$ws.Range("I132").Value = 3687.3333
$ws.Range("K132").Value = 11061.9999
$ws.Range("M132").Value = -8531.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2375396
$ws.Range("I7").Value = 3167102.8
$ws.Range("J7").Value = 275
$ws.Range("K7").Value = 3167102.8
$ws.Range("L7").Value = 275
$ws.Range("M7").Value = -3166989.8
$ws.Range("N7").Value = -501
$ws.Range("H12").Value = 1070.5
$ws.Range("J12").Value = 1416.3334
$ws.Range("L12").Value = 1416.3334
$ws.Range("N12").Value = -1752.3334
$ws.Range("H14").Value = 8777.666999999999
$ws.Range("J14").Value = 5549.5
$ws.Range("L14").Value = 5549.5
$ws.Range("N14").Value = -5893.5
$ws.Range("H16").Value = 249.2
$ws.Range("I16").Value = 161.5
$ws.Range("K16").Value = 161.5
$ws.Range("M16").Value = 8.5
$ws.Range("H105").Value = 5080.7856
$ws.Range("I105").Value = 3947.4285
$ws.Range("J105").Value = 6214.143
$ws.Range("K105").Value = 3947.4285
$ws.Range("L105").Value = 6214.143
$ws.Range("M105").Value = -2200.4285
$ws.Range("N105").Value = -9708.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5841.273
$ws.Range("J4").Value = 9333
$ws.Range("L4").Value = 9333
$ws.Range("N4").Value = -9557
$ws.Range("H12").Value = 29029.7
$ws.Range("J12").Value = 36249.625
$ws.Range("L12").Value = 36249.625
$ws.Range("N12").Value = -36589.625
$ws.Range("H99").Value = 6674.8887
$ws.Range("I99").Value = 5624.6924
$ws.Range("K99").Value = 5624.6924
$ws.Range("M99").Value = -4126.6924
$ws.Range("H126").Value = 6674.8887
$ws.Range("I126").Value = 5624.6924
$ws.Range("K126").Value = 16874.0772
$ws.Range("M126").Value = -14404.0772
$ws.Range("H132").Value = 7239
$ws.Range("I132").Value = 4081.8462
$ws.Range("J132").Value = 11799.333
$ws.Range("K132").Value = 12245.5386
$ws.Range("L132").Value = 35397.999
$ws.Range("M132").Value = -9715.5386
$ws.Range("N132").Value = -40457.999
$ws.Range("H134").Value = 3300.5
$ws.Range("I134").Value = 2910.6
$ws.Range("K134").Value = 8731.799999999999
$ws.Range("M134").Value = -6196.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1458.3334
$ws.Range("J131").Value = 1600
$ws.Range("L131").Value = 4800
$ws.Range("N131").Value = -14880
$ws.Range("H134").Value = 4552.125
$ws.Range("I134").Value = 4552.125
$ws.Range("K134").Value = 13656.375
$ws.Range("M134").Value = -8586.375
$ws.Range("H138").Value = 1250
$ws.Range("I138").Value = 1000
$ws.Range("J138").Value = 1500
$ws.Range("K138").Value = 3000
$ws.Range("L138").Value = 4500
$ws.Range("M138").Value = 2140
$ws.Range("N138").Value = -14780

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4750
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 4500
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -4645
$ws.Range("N41").Value = -5210
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 9899.666999999999
$ws.Range("J126").Value = 9899.666999999999
$ws.Range("L126").Value = 29699.001
$ws.Range("N126").Value = -34639.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5500
$ws.Range("I7").Value = 5500
$ws.Range("K7").Value = 5500
$ws.Range("M7").Value = -5388
$ws.Range("H22").Value = 2118.6667
$ws.Range("I22").Value = 1150
$ws.Range("K22").Value = 1150
$ws.Range("M22").Value = -855
$ws.Range("H27").Value = 2118.6667
$ws.Range("I27").Value = 1150
$ws.Range("K27").Value = 1150
$ws.Range("M27").Value = -1043
$ws.Range("H100").Value = 2249.75
$ws.Range("I100").Value = 1666.3334
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 1666.3334
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1125.3334
$ws.Range("N100").Value = -5082
$ws.Range("H116").Value = 247500
$ws.Range("J116").Value = 247500
$ws.Range("L116").Value = 247500
$ws.Range("N116").Value = -256678
$ws.Range("H122").Value = 7800
$ws.Range("I122").Value = 6600
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 19800
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -17350
$ws.Range("N122").Value = -31900
$ws.Range("H126").Value = 5500
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 4194
$ws.Range("I23").Value = 5130
$ws.Range("K23").Value = 5130
$ws.Range("M23").Value = -4901
$ws.Range("H57").Value = 95000
$ws.Range("J57").Value = 95000
$ws.Range("L57").Value = 95000
$ws.Range("N57").Value = -96508
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("H126").Value = 1450
$ws.Range("I126").Value = 1450
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4350
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1880
$ws.Range("N126").Value = $null
